$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.283732891082764
$ws.Range("B1").Value = 1.501825571060181
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.096885204315186
$ws.Range("E1").Value = 0.8715436458587646
